$wb = $excel.ActiveWorkbook

# --- Events sheet (sheet1) ---
$wsEvents = $wb.Worksheets.Item("Events")
$tbl1 = $wsEvents.ListObjects.Item("Table1")

# Rename existing B2 value
$wsEvents.Range("B2").Value = "Introduction to Microsoft Teams"

# Add Date, Start Time, End Time columns to Table1
$colDate = $tbl1.ListColumns.Add()
$colDate.Range.Item(1).Value = "Date"
$colStart = $tbl1.ListColumns.Add()
$colStart.Range.Item(1).Value = "Start Time"
$colEnd = $tbl1.ListColumns.Add()
$colEnd.Range.Item(1).Value = "End Time"

# Fill in data row
$wsEvents.Range("C2").Value = 43860
$wsEvents.Range("C2").NumberFormat = "mm-dd-yy"
$wsEvents.Range("D2").Value = 0.625
$wsEvents.Range("D2").NumberFormat = "h:mm"
$wsEvents.Range("E2").Value = 0.66666666666666663
$wsEvents.Range("E2").NumberFormat = "h:mm"

$wsEvents.Range("E3").Select()

$wsEvents.Columns.Item(1).ColumnWidth = 8.72
$wsEvents.Columns.Item(2).ColumnWidth = 31.17
$wsEvents.Range("C1:E1").ColumnWidth = 12.5

$tbl1.Resize($wsEvents.Range("A1:E3"))
